$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell C2 with the value "OK" (creates a new shared string entry)
$ws.Range("C2").Value = "OK"

# Update the active selection to C1 to match the target workbook state
$ws.Range("C1").Select()
